$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.966817792541496
$ws.Cells.Item(2, 3).Value = 0.893308662719607
$ws.Cells.Item(2, 4).Value = 0.584985586569476
$ws.Cells.Item(2, 5).Value = 0.296999999832692
$ws.Cells.Item(3, 2).Value = 0.919690806927313
$ws.Cells.Item(3, 3).Value = 0.7274036492185
$ws.Cells.Item(3, 4).Value = 0.211182613370219
$ws.Cells.Item(3, 5).Value = 0.0318949666208082
$ws.Cells.Item(4, 2).Value = 0.967698182118381
$ws.Cells.Item(4, 3).Value = 0.884062713988888
$ws.Cells.Item(4, 4).Value = 0.51440683255644
$ws.Cells.Item(4, 5).Value = 0.204280056564527
$ws.Cells.Item(5, 2).Value = 0.96149581434523
$ws.Cells.Item(5, 3).Value = 0.870539454458683
$ws.Cells.Item(5, 4).Value = 0.508864582772957
$ws.Cells.Item(5, 5).Value = 0.221856729679785
$ws.Cells.Item(6, 2).Value = 0.975960680252197
$ws.Cells.Item(6, 3).Value = 0.916588690371036
$ws.Cells.Item(6, 4).Value = 0.633780937303272
$ws.Cells.Item(6, 5).Value = 0.335309556593029
$ws.Cells.Item(7, 2).Value = 0.922424572474048
$ws.Cells.Item(7, 3).Value = 0.772863238082313
$ws.Cells.Item(7, 4).Value = 0.360351563726806
$ws.Cells.Item(7, 5).Value = 0.130819410340778
$ws.Cells.Item(8, 2).Value = 0.947296209112797
$ws.Cells.Item(8, 3).Value = 0.847810484166333
$ws.Cells.Item(8, 4).Value = 0.494485378656497
$ws.Cells.Item(8, 5).Value = 0.215296478222733
$ws.Cells.Item(9, 2).Value = 0.937765814921854
$ws.Cells.Item(9, 3).Value = 0.817919633035557
$ws.Cells.Item(9, 4).Value = 0.424406085295198
$ws.Cells.Item(9, 5).Value = 0.160883337220729
$ws.Cells.Item(10, 2).Value = 0.988227156408555
$ws.Cells.Item(10, 3).Value = 0.958443554017704
$ws.Cells.Item(10, 4).Value = 0.782134236066716
$ws.Cells.Item(10, 5).Value = 0.523181568917206
$ws.Cells.Item(11, 2).Value = 0.967620770835781
$ws.Cells.Item(11, 3).Value = 0.901215388289311
$ws.Cells.Item(11, 4).Value = 0.627650744781192
$ws.Cells.Item(11, 5).Value = 0.354631998553338
$ws.Cells.Item(12, 2).Value = 0.969048869019999
$ws.Cells.Item(12, 3).Value = 0.89743074638932
$ws.Cells.Item(12, 4).Value = 0.592405843268542
$ws.Cells.Item(12, 5).Value = 0.301973257468215
$ws.Cells.Item(13, 2).Value = 0.959259083655129
$ws.Cells.Item(13, 3).Value = 0.859802242608444
$ws.Cells.Item(13, 4).Value = 0.45949791657163
$ws.Cells.Item(13, 5).Value = 0.1626410563847
$ws.Cells.Item(14, 2).Value = 0.912219683163776
$ws.Cells.Item(14, 3).Value = 0.724249175069042
$ws.Cells.Item(14, 4).Value = 0.238062448920621
$ws.Cells.Item(14, 5).Value = 0.0488707424074325
$ws.Cells.Item(15, 2).Value = 0.958251568358937
$ws.Cells.Item(15, 3).Value = 0.875265260847593
$ws.Cells.Item(15, 4).Value = 0.564464818292818
$ws.Cells.Item(15, 5).Value = 0.290736339421869
$ws.Cells.Item(16, 2).Value = 0.946045661634544
$ws.Cells.Item(16, 3).Value = 0.836254446562586
$ws.Cells.Item(16, 4).Value = 0.463891749770764
$ws.Cells.Item(16, 5).Value = 0.199710949736978
$ws.Cells.Item(17, 2).Value = 0.90363865165016
$ws.Cells.Item(17, 3).Value = 0.734308986327845
$ws.Cells.Item(17, 4).Value = 0.282050474396961
$ws.Cells.Item(17, 5).Value = 0.0704404235784649
$ws.Cells.Item(18, 2).Value = 0.925036348492397
$ws.Cells.Item(18, 3).Value = 0.79149611071681
$ws.Cells.Item(18, 4).Value = 0.391796798922505
$ws.Cells.Item(18, 5).Value = 0.148366607425798
$ws.Cells.Item(19, 2).Value = 0.957381150247086
$ws.Cells.Item(19, 3).Value = 0.855842890735909
$ws.Cells.Item(19, 4).Value = 0.455796752262887
$ws.Cells.Item(19, 5).Value = 0.162838023869053
$ws.Cells.Item(20, 2).Value = 0.975928937846249
$ws.Cells.Item(20, 3).Value = 0.919277502372317
$ws.Cells.Item(20, 4).Value = 0.648883215002073
$ws.Cells.Item(20, 5).Value = 0.353123550451293
$ws.Cells.Item(21, 2).Value = 0.904399174743538
$ws.Cells.Item(21, 3).Value = 0.718967776890331
$ws.Cells.Item(21, 4).Value = 0.240193298194484
$ws.Cells.Item(21, 5).Value = 0.0486131668129673
$ws.Cells.Item(22, 2).Value = 0.955554705180873
$ws.Cells.Item(22, 3).Value = 0.874928348732178
$ws.Cells.Item(22, 4).Value = 0.583470561392379
$ws.Cells.Item(22, 5).Value = 0.318247605859147
$ws.Cells.Item(23, 2).Value = 0.932525596361294
$ws.Cells.Item(23, 3).Value = 0.810279176060404
$ws.Cells.Item(23, 4).Value = 0.427446499530605
$ws.Cells.Item(23, 5).Value = 0.174195954956821
$ws.Cells.Item(24, 2).Value = 0.947183440535729
$ws.Cells.Item(24, 3).Value = 0.844365613039134
$ws.Cells.Item(24, 4).Value = 0.497705718972695
$ws.Cells.Item(24, 5).Value = 0.243161072526148
$ws.Cells.Item(25, 2).Value = 0.902780968562746
$ws.Cells.Item(25, 3).Value = 0.729067526468771
$ws.Cells.Item(25, 4).Value = 0.275662643005539
$ws.Cells.Item(25, 5).Value = 0.0695438082440016
$ws.Cells.Item(26, 2).Value = 0.974508878224826
$ws.Cells.Item(26, 3).Value = 0.914972383100543
$ws.Cells.Item(26, 4).Value = 0.638218770811151
$ws.Cells.Item(26, 5).Value = 0.346199284268555
$ws.Cells.Item(27, 2).Value = 0.946736023881062
$ws.Cells.Item(27, 3).Value = 0.837583238786676
$ws.Cells.Item(27, 4).Value = 0.458310126767143
$ws.Cells.Item(27, 5).Value = 0.190738059135118
$ws.Cells.Item(28, 2).Value = 0.897331208873243
$ws.Cells.Item(28, 3).Value = 0.691923088574233
$ws.Cells.Item(28, 4).Value = 0.216828266125095
$ws.Cells.Item(28, 5).Value = 0.0471340124920385
$ws.Cells.Item(29, 2).Value = 0.91274005649211
$ws.Cells.Item(29, 3).Value = 0.754547632928866
$ws.Cells.Item(29, 4).Value = 0.309328477901705
$ws.Cells.Item(29, 5).Value = 0.0838624431102125
$ws.Cells.Item(30, 2).Value = 0.976070582671875
$ws.Cells.Item(30, 3).Value = 0.905950065437401
$ws.Cells.Item(30, 4).Value = 0.556329352719178
$ws.Cells.Item(30, 5).Value = 0.231833608409257
$ws.Cells.Item(31, 2).Value = 0.875704555382598
$ws.Cells.Item(31, 3).Value = 0.661156629633053
$ws.Cells.Item(31, 4).Value = 0.18554087292582
$ws.Cells.Item(31, 5).Value = 0.030740884870696
$ws.Cells.Item(32, 2).Value = 0.955525811468161
$ws.Cells.Item(32, 3).Value = 0.87448397227179
$ws.Cells.Item(32, 4).Value = 0.575624914173074
$ws.Cells.Item(32, 5).Value = 0.306448786346713
$ws.Cells.Item(33, 2).Value = 0.976193373324002
$ws.Cells.Item(33, 3).Value = 0.914741822758933
$ws.Cells.Item(33, 4).Value = 0.614657375837852
$ws.Cells.Item(33, 5).Value = 0.302394474824809
$ws.Cells.Item(34, 2).Value = 0.965835579071206
$ws.Cells.Item(34, 3).Value = 0.897904747656214
$ws.Cells.Item(34, 4).Value = 0.624212819204154
$ws.Cells.Item(34, 5).Value = 0.354229538715618
$ws.Cells.Item(35, 2).Value = 0.963593322408334
$ws.Cells.Item(35, 3).Value = 0.881072999822109
$ws.Cells.Item(35, 4).Value = 0.551679642770457
$ws.Cells.Item(35, 5).Value = 0.267289802945477
$ws.Cells.Item(36, 2).Value = 0.908648520864565
$ws.Cells.Item(36, 3).Value = 0.749997260704429
$ws.Cells.Item(36, 4).Value = 0.31061086318466
$ws.Cells.Item(36, 5).Value = 0.0862392811174852
$ws.Cells.Item(37, 2).Value = 0.951994434998097
$ws.Cells.Item(37, 3).Value = 0.853350063210115
$ws.Cells.Item(37, 4).Value = 0.496472862899741
$ws.Cells.Item(37, 5).Value = 0.219256641392608
$ws.Cells.Item(38, 2).Value = 0.958784174555761
$ws.Cells.Item(38, 3).Value = 0.866700786593034
$ws.Cells.Item(38, 4).Value = 0.502529074582536
$ws.Cells.Item(38, 5).Value = 0.207645400649768
$ws.Cells.Item(39, 2).Value = 0.944141220814485
$ws.Cells.Item(39, 3).Value = 0.806394795513822
$ws.Cells.Item(39, 4).Value = 0.339110724715792
$ws.Cells.Item(39, 5).Value = 0.0859356756796731
$ws.Cells.Item(40, 2).Value = 0.935226528278311
$ws.Cells.Item(40, 3).Value = 0.81873593979017
$ws.Cells.Item(40, 4).Value = 0.449249541656059
$ws.Cells.Item(40, 5).Value = 0.198796804292129
$ws.Cells.Item(41, 2).Value = 0.950403770486075
$ws.Cells.Item(41, 3).Value = 0.857333247582139
$ws.Cells.Item(41, 4).Value = 0.529505466967136
$ws.Cells.Item(41, 5).Value = 0.256941179329393
$ws.Cells.Item(42, 2).Value = 0.961456507874698
$ws.Cells.Item(42, 3).Value = 0.872809620910435
$ws.Cells.Item(42, 4).Value = 0.526657313337644
$ws.Cells.Item(42, 5).Value = 0.244132070491811
$ws.Cells.Item(43, 2).Value = 0.907353075411862
$ws.Cells.Item(43, 3).Value = 0.737785420397991
$ws.Cells.Item(43, 4).Value = 0.290389244830565
$ws.Cells.Item(43, 5).Value = 0.0796056363899287
$ws.Cells.Item(44, 2).Value = 0.957627685870134
$ws.Cells.Item(44, 3).Value = 0.859647989955396
$ws.Cells.Item(44, 4).Value = 0.486529442605217
$ws.Cells.Item(44, 5).Value = 0.205070127497966
$ws.Cells.Item(45, 2).Value = 0.972888464239402
$ws.Cells.Item(45, 3).Value = 0.911801389343149
$ws.Cells.Item(45, 4).Value = 0.639795127491623
$ws.Cells.Item(45, 5).Value = 0.353426469696419
$ws.Cells.Item(46, 2).Value = 0.904995459277873
$ws.Cells.Item(46, 3).Value = 0.715488052117734
$ws.Cells.Item(46, 4).Value = 0.228713125026319
$ws.Cells.Item(46, 5).Value = 0.0420365123887735
$ws.Cells.Item(47, 2).Value = 0.878069863206414
$ws.Cells.Item(47, 3).Value = 0.651428403703649
$ws.Cells.Item(47, 4).Value = 0.161812694473822
$ws.Cells.Item(47, 5).Value = 0.0221476380439832
$ws.Cells.Item(48, 2).Value = 0.940114490145502
$ws.Cells.Item(48, 3).Value = 0.828192713875771
$ws.Cells.Item(48, 4).Value = 0.467407800339511
$ws.Cells.Item(48, 5).Value = 0.214253701827007
$ws.Cells.Item(49, 2).Value = 0.923867619762656
$ws.Cells.Item(49, 3).Value = 0.760946181018884
$ws.Cells.Item(49, 4).Value = 0.289980966422097
$ws.Cells.Item(49, 5).Value = 0.0717323919720419
$ws.Cells.Item(50, 2).Value = 0.946818171036012
$ws.Cells.Item(50, 3).Value = 0.826258932428088
$ws.Cells.Item(50, 4).Value = 0.397656477588727
$ws.Cells.Item(50, 5).Value = 0.125206397620568
$ws.Cells.Item(51, 2).Value = 0.913713855881119
$ws.Cells.Item(51, 3).Value = 0.75627846775342
$ws.Cells.Item(51, 4).Value = 0.33939357519954
$ws.Cells.Item(51, 5).Value = 0.127829024311818
$ws.Cells.Item(52, 2).Value = 0.879443045320427
$ws.Cells.Item(52, 3).Value = 0.660719970715417
$ws.Cells.Item(52, 4).Value = 0.179995756591528
$ws.Cells.Item(52, 5).Value = 0.0292049872495772
